$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 42, pushing the existing rows 42-44 down to 44-46
$ws.Rows("42:43").Insert()

# Row 42: new "Especial" record (2022-02-17, serial 44615)
$ws.Range("A42").Value = 2
$ws.Range("B42").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C42").Value = "Coquimbo"
$ws.Range("D42").Value = 44615
$ws.Range("E42").Value = 4
$ws.Range("F42").Value = "Fruta"
$ws.Range("G42").Value = 100107
$ws.Range("H42").Value = "Otros"
$ws.Range("I42").Value = 100107011
$ws.Range("J42").Value = "Tuna"
$ws.Range("K42").Value = "Sin especificar"
$ws.Range("L42").Value = "Especial"
$ws.Range("M42").Value = 200
$ws.Range("N42").Value = 14000
$ws.Range("O42").Value = 15000
$ws.Range("P42").Value = 14500
$ws.Range("Q42").Value = "$/caja 18 kilos"
$ws.Range("R42").Value = "Provincia de Limarí"
$ws.Range("S42").Value = 806
$ws.Range("T42").Value = 18

# Row 43: new "Primera" record (2022-02-17, serial 44615)
$ws.Range("A43").Value = 2
$ws.Range("B43").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C43").Value = "Coquimbo"
$ws.Range("D43").Value = 44615
$ws.Range("E43").Value = 4
$ws.Range("F43").Value = "Fruta"
$ws.Range("G43").Value = 100107
$ws.Range("H43").Value = "Otros"
$ws.Range("I43").Value = 100107011
$ws.Range("J43").Value = "Tuna"
$ws.Range("K43").Value = "Sin especificar"
$ws.Range("L43").Value = "Primera"
$ws.Range("M43").Value = 400
$ws.Range("N43").Value = 12000
$ws.Range("O43").Value = 13000
$ws.Range("P43").Value = 12500
$ws.Range("Q43").Value = "$/caja 18 kilos"
$ws.Range("R43").Value = "Provincia de Limarí"
$ws.Range("S43").Value = 694
$ws.Range("T43").Value = 18
